$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H40").Value = 1800.3334
$ws.Range("J40").Value = 1700
$ws.Range("L40").Value = 1700
$ws.Range("N40").Value = -2050
$ws.Range("H55").Value = 218.75
$ws.Range("I55").Value = 100
$ws.Range("J55").Value = 235.71428
$ws.Range("K55").Value = 100
$ws.Range("L55").Value = 235.71428
$ws.Range("M55").Value = 114
$ws.Range("N55").Value = -663.71428
$ws.Range("H58").Value = 8627.416999999999
$ws.Range("H96").Value = 580
$ws.Range("I96").Value = 239
$ws.Range("K96").Value = 717
$ws.Range("M96").Value = 656
$ws.Range("H112").Value = 10205653
$ws.Range("J112").Value = 1603.9684
$ws.Range("L112").Value = 4811.9052
$ws.Range("N112").Value = -7027.9052
$ws.Range("H132").Value = 108276.92
$ws.Range("I132").Value = 130855.91
$ws.Range("K132").Value = 392567.73
$ws.Range("M132").Value = -390037.73
$ws.Range("H138").Value = 2448.5193
$ws.Range("I138").Value = 1298.2333
$ws.Range("J138").Value = 4017.0908
$ws.Range("K138").Value = 3894.699900000001
$ws.Range("L138").Value = 12051.2724
$ws.Range("M138").Value = 1245.300099999999
$ws.Range("N138").Value = -22331.2724

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H45").Value = 1535.8518
$ws.Range("I45").Value = 846.2857
$ws.Range("K45").Value = 846.2857
$ws.Range("M45").Value = -469.2857
$ws.Range("H122").Value = 3360.8333
$ws.Range("I122").Value = 1825
$ws.Range("J122").Value = 4896.6665
$ws.Range("K122").Value = 5475
$ws.Range("L122").Value = 14689.9995
$ws.Range("M122").Value = -3025
$ws.Range("N122").Value = -19589.9995

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 3055.932
$ws.Range("I31").Value = 1351.7368
$ws.Range("J31").Value = 4351.12
$ws.Range("K31").Value = 1351.7368
$ws.Range("L31").Value = 4351.12
$ws.Range("M31").Value = -1056.7368
$ws.Range("N31").Value = -4941.12
$ws.Range("H34").Value = 3055.932
$ws.Range("I34").Value = 1351.7368
$ws.Range("J34").Value = 4351.12
$ws.Range("K34").Value = 1351.7368
$ws.Range("L34").Value = 4351.12
$ws.Range("M34").Value = -1149.7368
$ws.Range("N34").Value = -4755.12
$ws.Range("H58").Value = 1782.3816
$ws.Range("I58").Value = 1555.2153
$ws.Range("J58").Value = 3124.7273
$ws.Range("K58").Value = 1555.2153
$ws.Range("L58").Value = 3124.7273
$ws.Range("M58").Value = -1352.2153
$ws.Range("N58").Value = -3530.7273
$ws.Range("H59").Value = 25056.182
$ws.Range("I59").Value = 6500
$ws.Range("J59").Value = 29179.777
$ws.Range("K59").Value = 6500
$ws.Range("L59").Value = 29179.777
$ws.Range("M59").Value = -5355
$ws.Range("N59").Value = -31469.777
$ws.Range("H94").Value = 1771.8572
$ws.Range("I94").Value = 945
$ws.Range("J94").Value = 2102.6
$ws.Range("K94").Value = 945
$ws.Range("L94").Value = 2102.6
$ws.Range("M94").Value = -494
$ws.Range("N94").Value = -3004.6
$ws.Range("H136").Value = 1782.3816
$ws.Range("I136").Value = 1555.2153
$ws.Range("J136").Value = 3124.7273
$ws.Range("K136").Value = 4665.6459
$ws.Range("L136").Value = 9374.1819
$ws.Range("M136").Value = -2115.6459
$ws.Range("N136").Value = -14474.1819

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H62").Value = 3000
$ws.Range("H65").Value = 3000
$ws.Range("H86").Value = 5089.5557
$ws.Range("I86").Value = 4776.923
$ws.Range("J86").Value = 5902.4
$ws.Range("K86").Value = 14330.769
$ws.Range("L86").Value = 17707.2
$ws.Range("M86").Value = -13144.769
$ws.Range("N86").Value = -20079.2
$ws.Range("H89").Value = 5089.5557
$ws.Range("I89").Value = 4776.923
$ws.Range("J89").Value = 5902.4
$ws.Range("K89").Value = 42992.307
$ws.Range("L89").Value = 53121.6
$ws.Range("M89").Value = -37064.307
$ws.Range("N89").Value = -64977.6
$ws.Range("H131").Value = 17858266
$ws.Range("I131").Value = 125002690
$ws.Range("J131").Value = 862.2917
$ws.Range("K131").Value = 375008070
$ws.Range("L131").Value = 2586.8751
$ws.Range("M131").Value = -375003030
$ws.Range("N131").Value = -12666.8751

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 6710.6665
$ws.Range("I70").Value = 5979.1304
$ws.Range("K70").Value = 5979.1304
$ws.Range("M70").Value = -5709.1304
$ws.Range("H73").Value = 6710.6665
$ws.Range("I73").Value = 5979.1304
$ws.Range("K73").Value = 5979.1304
$ws.Range("M73").Value = -5043.1304
$ws.Range("H80").Value = 25002472
$ws.Range("I80").Value = 41668950
$ws.Range("J80").Value = 2750
$ws.Range("K80").Value = 41668950
$ws.Range("L80").Value = 2750
$ws.Range("M80").Value = -41667952
$ws.Range("N80").Value = -4746
$ws.Range("H83").Value = 25002472
$ws.Range("I83").Value = 41668950
$ws.Range("J83").Value = 2750
$ws.Range("K83").Value = 208344750
$ws.Range("L83").Value = 13750
$ws.Range("M83").Value = -208339758
$ws.Range("N83").Value = -23734
$ws.Range("H122").Value = 2672.0881
$ws.Range("I122").Value = 2324
$ws.Range("J122").Value = 3113
$ws.Range("K122").Value = 6972
$ws.Range("L122").Value = 9339
$ws.Range("M122").Value = -4522
$ws.Range("N122").Value = -14239
$ws.Range("H136").Value = 11339.577
$ws.Range("J136").Value = 11769.56
$ws.Range("L136").Value = 35308.68
$ws.Range("N136").Value = -40408.68

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H40").Value = 5793.4814
$ws.Range("I40").Value = 5318.8945
$ws.Range("J40").Value = 6920.625
$ws.Range("K40").Value = 5318.8945
$ws.Range("L40").Value = 6920.625
$ws.Range("M40").Value = -5182.8945
$ws.Range("N40").Value = -7192.625
$ws.Range("H69").Value = 356721
$ws.Range("J69").Value = 356721
$ws.Range("L69").Value = 356721
$ws.Range("N69").Value = -358343
$ws.Range("H72").Value = 356721
$ws.Range("J72").Value = 356721
$ws.Range("L72").Value = 1070163
$ws.Range("N72").Value = -1078275
$ws.Range("H122").Value = 3408.9644
$ws.Range("I122").Value = 2691.7083
$ws.Range("J122").Value = 7712.5
$ws.Range("K122").Value = 8075.124899999999
$ws.Range("L122").Value = 23137.5
$ws.Range("M122").Value = -5625.124899999999
$ws.Range("N122").Value = -28037.5

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H47").Value = 26250
$ws.Range("J47").Value = 26250
$ws.Range("L47").Value = 26250
$ws.Range("N47").Value = -27394
$ws.Range("H113").Value = 209.75
$ws.Range("I113").Value = 209.75
$ws.Range("K113").Value = 629.25
$ws.Range("M113").Value = 1540.75
$ws.Range("H132").Value = 5465809
$ws.Range("I132").Value = 492.48718
$ws.Range("J132").Value = 15154325
$ws.Range("K132").Value = 1477.46154
$ws.Range("L132").Value = 45462975
$ws.Range("M132").Value = 1052.53846
$ws.Range("N132").Value = -45468035
